# Wed, Jun 24, 2020  2:06:04 PM
#
# 1) The table on slide 16 ("Google Shape;213;p29") switches from the
#    deck's custom "Table_0" style to the built-in PowerPoint table
#    style {8DB347C9-1D4D-4FE5-8549-D3EFC81DCF79}.
# 2) The theme colour scheme used by the deck (ppt/theme/theme1.xml,
#    attached to the one-and-only Slide Master) is swapped from the
#    "Integral" palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{8DB347C9-1D4D-4FE5-8549-D3EFC81DCF79}")

# --- 2. Swap the theme palette from Integral -> Office ----------------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

function ToRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colours, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Item(1).RGB  = ToRgb 0x00 0x00 0x00   # dk1      000000
$colors.Item(2).RGB  = ToRgb 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Item(3).RGB  = ToRgb 0x44 0x54 0x6A   # dk2      44546A
$colors.Item(4).RGB  = ToRgb 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Item(5).RGB  = ToRgb 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Item(6).RGB  = ToRgb 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Item(7).RGB  = ToRgb 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Item(8).RGB  = ToRgb 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Item(9).RGB  = ToRgb 0x44 0x72 0xC4   # accent5  4472C4
$colors.Item(10).RGB = ToRgb 0x70 0xAD 0x47   # accent6  70AD47
$colors.Item(11).RGB = ToRgb 0x05 0x63 0xC1   # hlink    0563C1
$colors.Item(12).RGB = ToRgb 0x95 0x4F 0x72   # folHlink 954F72
